$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cell H1 (bold, centered, thin-bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the new I and J data columns for rows 2-73.
$iVals = @(6,7,7,8,8,7,7,5,7,6,6,7,7,6,7,8,7,7,8,8,7,7,7,8,7,7,9,7,7,7,7,7,7,5,7,7,6,8,8,10,8,8,7,8,6,8,8,8,8,7,6,7,7,8,9,7,8,6,7,6,7,8,7,7,8,8,7,7,6,6,6,6)
$jVals = @(7,8,7,8,8,7,7,5,7,6,6,7,7,6,7,8,7,7,8,8,8,8,7,8,7,7,9,7,7,7,7,7,7,5,8,7,6,8,8,10,8,8,7,8,6,8,8,9,9,8,6,8,8,8,9,7,8,7,7,6,7,8,7,7,8,8,7,7,6,6,6,6)

for ($r = 2; $r -le 73; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
